$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows to write, in order, starting at row 2
# Columns: A (Região), B (Variável), C (Valor), D (Colocação, text with º or blank)
$rows = @(
  @("Paraíba", "Diferença 2022-2013", 0.04500000000000004, "1º"),
  @("Rio de Janeiro", "Diferença 2022-2013", 0.02200000000000002, "2º"),
  @("Rio Grande do Norte", "Diferença 2022-2013", 0.02100000000000002, "3º"),
  @("Amapá", "Diferença 2022-2013", 0.01800000000000002, "4º"),
  @("Mato Grosso do Sul", "Diferença 2022-2013", 0.006000000000000005, "5º"),
  @("Tocantins", "Diferença 2022-2013", -0.008000000000000007, "10º"),
  @("Goiás", "Diferença 2022-2013", -0.009000000000000008, "11º"),
  @("Roraima", "Diferença 2022-2013", -0.009000000000000008, "12º"),
  @("Santa Catarina", "Diferença 2022-2013", -0.009000000000000008, "13º"),
  @("Rio Grande do Sul", "Diferença 2022-2013", -0.009999999999999953, "14º"),
  @("São Paulo", "Diferença 2022-2013", -0.01000000000000001, "15º"),
  @("Alagoas", "Diferença 2022-2013", -0.01200000000000001, "16º"),
  @("Piauí", "Diferença 2022-2013", -0.01700000000000002, "17º"),
  @("Ceará", "Diferença 2022-2013", -0.01900000000000002, "18º"),
  @("Pará", "Diferença 2022-2013", -0.02300000000000002, "19º"),
  @("Bahia", "Diferença 2022-2013", -0.02400000000000002, "20º"),
  @("Minas Gerais", "Diferença 2022-2013", -0.02599999999999997, "21º"),
  @("Rondônia", "Diferença 2022-2013", -0.02799999999999997, "22º"),
  @("Acre", "Diferença 2022-2013", -0.03000000000000003, "23º"),
  @("Distrito Federal", "Diferença 2022-2013", -0.04299999999999993, "24º"),
  @("Sergipe", "Diferença 2022-2013", -0.04399999999999993, "25º"),
  @("Pernambuco", "Diferença 2022-2013", -0.05099999999999993, "26º"),
  @("Amazonas", "Diferença 2022-2013", -0.07199999999999995, "27º"),
  @("Nordeste", "Diferença 2022-2013", -0.02100000000000002, $null),
  @("Brasil", "Diferença 2022-2013", -0.01400000000000001, $null)
)

$r = 2
foreach ($row in $rows) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  if ($null -ne $row[3]) {
    $ws.Cells.Item($r, 4).Value = $row[3]
  } else {
    $ws.Cells.Item($r, 4).Value = $null
  }
  $r = $r + 1
}
